# Updates need_to_buy.xlsx: refreshed forecast data (rows 2-15, columns A-F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45965, 8470.38170862993, 8034.57328943618, 9364, 3881.79105219784, 106.348514234751),
    @(45966, 8640.62761036021, 8292.16652585553, 9004, 4671.39971057666, 164.981926518008),
    @(45967, 8669.99553555781, 8482.28337007052, 9004, 4675.08988562196, 173.057218987187),
    @(45968, 8350.76523502114, 7595.99780372264, 9004, 4498.66269890353, 128.777520942757),
    @(45969, 2936.84492903534, 5263.11485391265, 9004, 4119.30896387416, 15.7676590744506),
    @(45970, 2982.76489981526, 5152.80239426028, 9004, 4297.46132010564, 18.5943214319135),
    @(45971, 9273.47168962819, 8492.35768589662, 9004, 5243.44716672023, 197.158535525702),
    @(45972, 9273.47299076555, 8703.20822710807, 9004, 5243.45092906603, 205.944131507254),
    @(45973, 9273.47299076555, 8359.22333081911, 9004, 5243.45092906603, 191.611427495214),
    @(45974, 9273.47299076555, 8173.19352522393, 9004, 5243.45092906603, 183.860185595415),
    @(45975, 9273.47299076555, 7037.40376858555, 9004, 5243.45092906603, 136.535612402149),
    @(45976, 3551.20396235217, 4390.17868219587, 9004, 4928.56324744655, 13.1142470684342),
    @(45977, 3412.56435984501, 4300.49389936102, 9004, 4921.0914342312,  9.06605556634258),
    @(45978, 10116.508538002,  7836.01442000187, 9004, 5829.29599708779, 194.221267378736)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
